$wb = $excel.ActiveWorkbook

# Add the new worksheet for the "soglia diodo discendente" data
$ws = $wb.Worksheets.Add()
$ws.Name = "soglia diodo discendente"

# Move it to the end of the tab strip (after "soglia diodo crescenti")
$wsMove = $wb.Worksheets.Item("soglia diodo discendente")
$lastSheet = $wb.Worksheets.Item("soglia diodo crescenti")
$wsMove.Move($null, $lastSheet)

# Re-fetch the worksheet reference (position/identity may have changed after
# the move) and make it the active sheet, as in the source workbook.
$ws = $wb.Worksheets.Item("soglia diodo discendente")
$ws.Activate()

# Header row
$ws.Range("A1").Value = "V"
$ws.Range("B1").Value = "I(uA)"
$ws.Range("C1").Value = "NOTA: presi in ordine decrescente"

# Data rows 2-26 (V in column A, I(uA) in column B)
$data = New-Object 'object[,]' 25,2
$data[0,0] = 0.99399999999999999
$data[0,1] = 83950
$data[1,0] = 0.98399999999999999
$data[1,1] = 80623
$data[2,0] = 0.93899999999999995
$data[2,1] = 65630
$data[3,0] = 0.91400000000000003
$data[3,1] = 57443
$data[4,0] = 0.89400000000000002
$data[4,1] = 51095
$data[5,0] = 0.88400000000000001
$data[5,1] = 47950
$data[6,0] = 0.86899999999999999
$data[6,1] = 43340
$data[7,0] = 0.83499999999999996
$data[7,1] = 33040
$data[8,0] = 0.80500000000000005
$data[8,1] = 24780
$data[9,0] = 0.79400000000000004
$data[9,1] = 21925
$data[10,0] = 0.78000000000000003
$data[10,1] = 18460
$data[11,0] = 0.77000000000000002
$data[11,1] = 16165
$data[12,0] = 0.76000000000000001
$data[12,1] = 14020
$data[13,0] = 0.74399999999999999
$data[13,1] = 11045
$data[14,0] = 0.72999999999999998
$data[14,1] = 8509
$data[15,0] = 0.71399999999999997
$data[15,1] = 6312
$data[16,0] = 0.68899999999999995
$data[16,1] = 3580
$data[17,0] = 0.64300000000000002
$data[17,1] = 1080
$data[18,0] = 0.61799999999999999
$data[18,1] = 190
$data[19,0] = 0.60699999999999998
$data[19,1] = 161
$data[20,0] = 0.59199999999999997
$data[20,1] = 127
$data[21,0] = 0.52100000000000002
$data[21,1] = 31.379999999999999
$data[22,0] = 0.46899999999999997
$data[22,1] = 8.1300000000000008
$data[23,0] = 0.38700000000000001
$data[23,1] = 0.77000000000000002
$data[24,0] = 0.224
$data[24,1] = 0.01
$ws.Range("A2:B26").Value = $data

# Number formats matching the source layout:
#  - row 2 keeps the default "General" format (as in the source)
#  - A3:A26 -> "0.000"
#  - B3:B22 -> "0" (integer uA readings)
#  - B23:B26 -> "0.00" (sub-1 uA readings)
$ws.Range("A3:A26").NumberFormat = "0.000"
$ws.Range("B3:B22").NumberFormat = "0"
$ws.Range("B23:B26").NumberFormat = "0.00"

$ws.Range("A3").Select()
